$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AMS_data_summary")
$v = $ws.Cells.Item(1,1).Value()
Write-Output ("cell1,1=" + $v)
$v2 = $ws.Range("A1").Value()
Write-Output ("A1=" + $v2)
